# Field Mapping code Commit
# Adds a new test row (row 18) to the "Login" sheet with a hyperlinked
# e-mail address in column B, matching the existing pattern used for
# e-mail-looking values elsewhere in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# New data row
$ws.Range("A18").Value = "Field_Mapping_TC001"
$ws.Range("B18").Value = "karthirko29@gmail.com"
$ws.Range("C18").Value = "Password1"
$ws.Range("D18").Value = "Login successful"

# Turn the e-mail address in B18 into a mailto hyperlink, and apply the
# workbook's built-in Hyperlink style to match the other hyperlinked cells.
[void]$ws.Hyperlinks.Add($ws.Range("B18"), "mailto:karthirko29@gmail.com")
$ws.Range("B18").Style = "Hyperlink"

# Move the active selection to B20, as in the final workbook state.
[void]$ws.Range("B20").Select()
